# Automatische test-sync: 2025-06-18 15:30:10
# Append two new "Vragen over samenwerking" rows to the Logs sheet and
# bump the "Overig" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$subject   = "Vragen over samenwerking"
$sender    = "mailmind.test@zohomail.eu"
$content   = "Kunnen we samenwerken aan een nieuw project?"
$category  = "Overig"
$timestamp = "2025-06-18 15:00:11"
$answered  = "Nee"

$newRows = @(15, 16)
foreach ($r in $newRows) {
    $logs.Cells.Item($r, 1).Value = $subject
    $logs.Cells.Item($r, 2).Value = $sender
    $logs.Cells.Item($r, 3).Value = $content
    $logs.Cells.Item($r, 4).Value = $category
    $logs.Cells.Item($r, 6).Value = $timestamp
    $logs.Cells.Item($r, 7).Value = $answered
}

# Update the "Overig" count on the Dashboard sheet (was 5, now 7).
$dashboard.Range("B2").Value = 7

# Extend the two conditional-formatting rules on the Logs sheet so they
# keep covering the newly appended rows (D2:D14 -> D2:D16, G2:G14 -> G2:G16).
$catConditions = $logs.Range("D2:D14").FormatConditions
for ($i = 1; $i -le $catConditions.Count; $i++) {
    $catConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D16"))
}

$answeredConditions = $logs.Range("G2:G14").FormatConditions
for ($i = 1; $i -le $answeredConditions.Count; $i++) {
    $answeredConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G16"))
}
